$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells remain text so numeric-looking strings (dates, multi-dot numbers,
# leading/trailing zeros) are preserved exactly as authored, matching the source data.
$cellUpdates = @{
    "D2" = "27.699.39"
    "E2" = "  +1.57%  "
    "D3" = "1.866.18"
    "E3" = "  +0.54%  "
    "D4" = "1.005"
    "E4" = "  +0.39%  "
    "D5" = "331.37"
    "E5" = "  +3.18%  "
    "E6" = "  +0.32%  "
    "D7" = "0.4673"
    "E7" = "  +3.97%  "
    "D8" = "0.3932"
    "E8" = "  +2.08%  "
    "D9" = "47.49"
    "E9" = "  -0.51%  "
    "D10" = "0.08015"
    "E10" = "  +1.56%  "
    "E11" = "  +0.31%  "
    "D12" = "21.73"
    "E12" = "  +1.87%  "
    "D13" = "1.892.52"
    "E13" = "  +2.36%  "
    "D14" = "5.925"
    "E14" = "  +0.87%  "
    "D15" = "7.117"
    "E15" = "  -0.62%  "
    "D16" = "1.004"
    "E16" = "  +0.39%  "
    "D17" = "0.00001047"
    "E17" = "  +1.66%  "
    "D18" = "86.57"
    "E18" = "  +1.09%  "
    "D19" = "0.06616"
    "E19" = "  +1.40%  "
    "D20" = "17.12"
    "E20" = "  +1.15%  "
    "E21" = "  +0.30%  "
    "D22" = "27.714.09"
    "E22" = "  +1.61%  "
    "D23" = "5.480"
    "E23" = "  -0.31%  "
    "D24" = "10.97"
    "E24" = "  +1.85%  "
    "E25" = "  +2.21%  "
    "D26" = "2.109.54"
    "E26" = "  +1.81%  "
    "D27" = "159.06"
    "E27" = "  +4.87%  "
    "D28" = "20.14"
    "E28" = "  +2.59%  "
    "E29" = "  +1.39%  "
    "D30" = "5.550"
    "E30" = "  +1.74%  "
    "D31" = "122.70"
    "E31" = "  +2.00%  "
    "D32" = "0.9662"
    "E32" = "  +3.12%  "
    "D33" = "0.09469"
    "E33" = "  +2.20%  "
    "E34" = "  -1.31%  "
    "D35" = "3.598"
    "E35" = "  +0.79%  "
    "D36" = "5.306"
    "E36" = "  +0.28%  "
    "D37" = "0.02256"
    "E37" = "  +1.44%  "
    "D38" = "0.06059"
    "E38" = "  +1.31%  "
    "E39" = "  +2.19%  "
    "D40" = "8.124"
    "E40" = "  -2.01%  "
    "E41" = "  +0.26%  "
    "D42" = "0.5972"
    "D43" = "0.1891"
    "E43" = "  +0.55%  "
    "D44" = "10.24"
    "E44" = "  +1.39%  "
    "D45" = "1.266"
    "E45" = "  +1.00%  "
    "D46" = "0.5696"
    "E46" = "  +1.34%  "
    "D47" = "12.15"
    "E47" = "  +2.07%  "
    "D48" = "3.390"
    "E48" = "  +1.12%  "
    "D49" = "1.931"
    "E49" = "  +0.78%  "
    "D50" = "0.06843"
    "E50" = "  +0.60%  "
    "E51" = "  +5.32%  "
}

foreach ($cellRef in $cellUpdates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $cellUpdates[$cellRef]
}
